# Apply "corrected most names to the official names from website" edit
# - Change District "Madhugiri" -> "Tumakuru (Tumkur)" in column G for the affected rows
# - Clear the stray empty F cells (F15, F28, F41, F48) so they are no longer present

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update District column (G) values from Madhugiri to Tumakuru (Tumkur)
$rowsToUpdate = @(12, 16, 21, 26, 27, 33, 49)
foreach ($r in $rowsToUpdate) {
    $ws.Range("G$r").Value = "Tumakuru (Tumkur)"
}

# Clear stray empty inline-string cells in column F
$rowsToClear = @(15, 28, 41, 48)
foreach ($r in $rowsToClear) {
    $ws.Range("F$r").ClearContents()
}
